$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 294
$ws1.Range("F4").Value = 3609
$ws1.Range("F5").Value = 2222
$ws1.Range("F7").Value = 4
$ws1.Range("F9").Value = 82
$ws1.Range("F11").Value = 1330
$ws1.Range("F13").Value = 1944

# Sheet "全部类型" (sheet4): update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 294
$ws4.Range("F4").Value = 3609
$ws4.Range("F5").Value = 2222
$ws4.Range("F7").Value = 4
$ws4.Range("F10").Value = 82
$ws4.Range("F14").Value = 1330
$ws4.Range("F16").Value = 1944
